$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two row re-orderings)
# Force text format on target cells first so numeric-looking strings
# (e.g. "0.999", "1.69") are preserved as text, matching the source inlineStr cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.742.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.697.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.27%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.79'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.50'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.06'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.403'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '30.24'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000202'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +10.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.181.71'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.614.84'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.682.83'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.71'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.88'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '359.80'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.56'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.10%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000106'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +12.32%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.69'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.70%  '
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'SuiNetwork'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.63'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.41%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.24'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '531.56'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.19%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.65'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.45'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.432'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.79'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '162.95'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.31%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.35%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '169.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.65'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.10%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0613'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.42'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.65%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.50%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0266'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.36%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.658'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.07'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +8.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0983'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.20%  '
